$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before row 944, pushing existing rows 944:1000
# down to 949:1005 (and the used range / dimension grows to R1005).
$ws.Range("A944:A948").EntireRow.Insert()

# Populate the 5 newly inserted rows (944-948) with the new weekly data.
$rows = @(944, 945, 946, 947, 948)
$D = @(45041, 45041, 45041, 45041, 45041)
$H = @("Copenhague", "Crespo record", "Crespo record", "Crespo record", "Morada(o)")
$J = @(200, 800, 1000, 800, 250)
$K = @(1400, 1400, 1400, 1300, 1400)
$L = @(1400, 1400, 1400, 1300, 1400)
$M = @(1400, 1400, 1400, 1300, 1400)
$O = @("Provincia de Cautín", "Provincia de Cautín", "Región Metropolitana", "Región de O'Higgins", "Provincia de Cautín")
$P = @(1400, 1400, 1400, 1300, 1400)

for ($i = 0; $i -lt 5; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 4).Value = $D[$i]
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = 100112006
    $ws.Cells.Item($r, 7).Value = "Repollo"
    $ws.Cells.Item($r, 8).Value = $H[$i]
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $J[$i]
    $ws.Cells.Item($r, 11).Value = $K[$i]
    $ws.Cells.Item($r, 12).Value = $L[$i]
    $ws.Cells.Item($r, 13).Value = $M[$i]
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = $O[$i]
    $ws.Cells.Item($r, 16).Value = $P[$i]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
